# Auto-generated edit script to update JudgeBotOpinions (C) and JudgeBotFunctionCalls (D) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 'MSG: None

MSG: The decision about which movie to show on Friday was not made.
'
$ws.Cells.Item(2, 4).Value = 'no_decision, '
$ws.Cells.Item(3, 3).Value = 'MSG: None

MSG: The decision-making process did not yield a concrete choice for Friday''s movie.
'
$ws.Cells.Item(3, 4).Value = 'no_decision, '
$ws.Cells.Item(4, 3).Value = 'MSG: None

MSG: The decision has been recorded as no decision regarding Friday''s movie.
'
$ws.Cells.Item(4, 4).Value = 'no_decision, '
$ws.Cells.Item(5, 3).Value = 'MSG: None

MSG: The decision has been recorded as no decision was made about the movie to be shown on Friday.
'
$ws.Cells.Item(5, 4).Value = 'no_decision, '
$ws.Cells.Item(6, 3).Value = 'MSG: None

MSG: The decision has been recorded with no movie selected for Friday.
'
$ws.Cells.Item(6, 4).Value = 'no_decision, '
$ws.Cells.Item(7, 3).Value = 'MSG: None

MSG: The decision has been recorded: "Oppenheimer" will be the movie shown on Friday.
'
$ws.Cells.Item(8, 3).Value = 'MSG: None

MSG: The decision process concluded without agreeing on a specific movie for Friday, so no movie will be acquired.
'
$ws.Cells.Item(8, 4).Value = 'no_decision, '
$ws.Cells.Item(9, 3).Value = 'MSG: None

MSG: I have successfully recorded the decision to acquire the rights for "Barbie" for the movie to be shown on Friday.
'
$ws.Cells.Item(10, 3).Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Oppenheimer."
'
$ws.Cells.Item(11, 3).Value = 'MSG: None

MSG: The decision reflects that there was no agreement on a movie to show on Friday. The conversation did not lead to a definitive choice, so the outcome is recorded accordingly.
'
$ws.Cells.Item(11, 4).Value = 'no_decision, '
$ws.Cells.Item(12, 3).Value = 'MSG: None

MSG: The decision has been recorded, and no movie was selected for showing on Friday.
'
$ws.Cells.Item(12, 4).Value = 'no_decision, '
$ws.Cells.Item(13, 3).Value = 'MSG: None

MSG: The decision regarding which movie to show on Friday has resulted in no clear consensus.
'
$ws.Cells.Item(13, 4).Value = 'no_decision, '
$ws.Cells.Item(14, 3).Value = 'MSG: None

MSG: I have recorded the decision to acquire the rights to "Barbie."
'
$ws.Cells.Item(15, 3).Value = 'MSG: None

MSG: The decision has been recorded as no decision being made regarding Friday''s movie.
'
$ws.Cells.Item(15, 4).Value = 'no_decision, '
$ws.Cells.Item(16, 3).Value = 'MSG: None

MSG: The decision process concluded without a clear choice for Friday''s movie, resulting in no decision being made.
'
$ws.Cells.Item(16, 4).Value = 'no_decision, '
$ws.Cells.Item(17, 3).Value = 'MSG: None

MSG: The decision has been noted as "no decision" regarding the movie to be shown on Friday.
'
$ws.Cells.Item(17, 4).Value = 'no_decision, '
$ws.Cells.Item(18, 3).Value = 'MSG: None

MSG: I have recorded the decision to acquire rights for both movies.
'
$ws.Cells.Item(18, 4).Value = 'both_movies, '
$ws.Cells.Item(19, 3).Value = 'MSG: None

MSG: I have recorded the decision to acquire the rights for "Barbie" as the movie to be shown on Friday.
'
$ws.Cells.Item(20, 3).Value = 'MSG: None

MSG: The decision has been registered as no decision regarding the movie for Friday was made.
'
$ws.Cells.Item(20, 4).Value = 'no_decision, '
$ws.Cells.Item(21, 3).Value = 'MSG: None

MSG: The decision to acquire the rights for "Barbie" has been successfully recorded.
'
$ws.Cells.Item(22, 3).Value = 'MSG: None

MSG: The decision has been recorded as no movie selected for Friday.
'
$ws.Cells.Item(22, 4).Value = 'no_decision, '
$ws.Cells.Item(23, 3).Value = 'MSG: None

MSG: The decision to acquire the rights for "Oppenheimer" has been successfully recorded.
'
$ws.Cells.Item(24, 3).Value = 'MSG: None

MSG: The decision has been recorded, indicating that no agreement was reached regarding the movie for Friday.
'
$ws.Cells.Item(24, 4).Value = 'no_decision, '
$ws.Cells.Item(25, 3).Value = 'MSG: None

MSG: The rights for both movies have been successfully acquired.
'
$ws.Cells.Item(25, 4).Value = 'both_movies, '
$ws.Cells.Item(26, 3).Value = 'MSG: None

MSG: The decision about the movie for Friday was not made.
'
$ws.Cells.Item(26, 4).Value = 'no_decision, '
$ws.Cells.Item(27, 3).Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Barbie."
'
$ws.Cells.Item(28, 3).Value = 'MSG: None

MSG: The decision-making process concluded without selecting a movie for Friday.
'
$ws.Cells.Item(28, 4).Value = 'no_decision, '
$ws.Cells.Item(29, 3).Value = 'MSG: None

MSG: The decision regarding the movie to be shown on Friday could not be made.
'
$ws.Cells.Item(29, 4).Value = 'no_decision, '
$ws.Cells.Item(30, 3).Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Barbie."
'
$ws.Cells.Item(31, 3).Value = 'MSG: None

MSG: The decision has been recorded, and no movie was chosen for Friday.
'
$ws.Cells.Item(31, 4).Value = 'no_decision, '
$ws.Cells.Item(32, 3).Value = 'MSG: None

MSG: The decision has been made to acquire the rights for both movies.
'
$ws.Cells.Item(32, 4).Value = 'both_movies, '
$ws.Cells.Item(33, 3).Value = 'MSG: None

MSG: The committee did not make a decision about which movie to show on Friday.
'
$ws.Cells.Item(33, 4).Value = 'no_decision, '
$ws.Cells.Item(34, 3).Value = 'MSG: None

MSG: The decision was successfully recorded to acquire the rights for "Barbie." The movie will be shown on Friday.
'
$ws.Cells.Item(35, 3).Value = 'MSG: None

MSG: The decision process has concluded with no agreement on a movie for Friday.
'
$ws.Cells.Item(35, 4).Value = 'no_decision, '
$ws.Cells.Item(36, 3).Value = 'MSG: None

MSG: The decision has been recorded as no consensus reached for Friday''s movie.
'
$ws.Cells.Item(36, 4).Value = 'no_decision, '
$ws.Cells.Item(37, 3).Value = 'MSG: None

MSG: The rights to both movies will be acquired as per the committee''s decision.
'
$ws.Cells.Item(37, 4).Value = 'both_movies, '
$ws.Cells.Item(38, 3).Value = 'MSG: None

MSG: The decision-making conversation did not result in a consensus on which movie to show on Friday, leading to no movie being selected.
'
$ws.Cells.Item(38, 4).Value = 'no_decision, '
$ws.Cells.Item(39, 3).Value = 'MSG: None

MSG: No decision was made regarding which movie to show on Friday.
'
$ws.Cells.Item(39, 4).Value = 'no_decision, '
$ws.Cells.Item(40, 3).Value = 'MSG: None

MSG: The movie "Barbie" has been selected for acquisition.
'
$ws.Cells.Item(41, 3).Value = 'MSG: None

MSG: The rights to both movies have been successfully acquired.
'
$ws.Cells.Item(41, 4).Value = 'both_movies, '
$ws.Cells.Item(42, 3).Value = 'MSG: None

MSG: The decision has been recorded as no agreement on the movie for Friday.
'
$ws.Cells.Item(43, 3).Value = 'MSG: None

MSG: The decision has been recorded with no movie selected for Friday.
'
$ws.Cells.Item(43, 4).Value = 'no_decision, '
$ws.Cells.Item(44, 3).Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Oppenheimer."
'
$ws.Cells.Item(45, 3).Value = 'MSG: None

MSG: The rights to both movies have been acquired successfully.
'
$ws.Cells.Item(45, 4).Value = 'both_movies, '
$ws.Cells.Item(46, 3).Value = 'MSG: None

MSG: The decision has been recorded, and the rights for "Oppenheimer" will be acquired.
'
$ws.Cells.Item(47, 3).Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Barbie."
'
$ws.Cells.Item(48, 3).Value = 'MSG: None

MSG: The decision to acquire the rights to "Barbie" has been made.
'
$ws.Cells.Item(49, 3).Value = 'MSG: None

MSG: The decision has been recorded as no agreement was reached about which movie to show on Friday.
'
$ws.Cells.Item(49, 4).Value = 'no_decision, '
$ws.Cells.Item(50, 3).Value = 'MSG: None

MSG: The decision process has concluded without reaching an agreement on which movie to show on Friday.
'
$ws.Cells.Item(50, 4).Value = 'no_decision, '
$ws.Cells.Item(51, 3).Value = 'MSG: None

MSG: The decision regarding the movie to show on Friday was not made, resulting in no definitive choice.
'
$ws.Cells.Item(51, 4).Value = 'no_decision, '
$ws.Cells.Item(52, 3).Value = 'MSG: None

MSG: The decision has been recorded as no decision regarding the movie for Friday.
'
$ws.Cells.Item(52, 4).Value = 'no_decision, '
$ws.Cells.Item(53, 3).Value = 'MSG: None

MSG: The decision has been recorded as no decision regarding the movie for Friday.
'
$ws.Cells.Item(53, 4).Value = 'no_decision, '
$ws.Cells.Item(54, 3).Value = 'MSG: None

MSG: The decision has been recorded as "no decision," indicating that the committee has not finalized their selection for the movie to be shown on Friday.
'
$ws.Cells.Item(54, 4).Value = 'no_decision, '
$ws.Cells.Item(55, 3).Value = 'MSG: None

MSG: The decision about Friday''s movie has been recorded as "no decision."
'
$ws.Cells.Item(55, 4).Value = 'no_decision, '
$ws.Cells.Item(56, 3).Value = 'MSG: None

MSG: It appears that there was no decision made regarding which movie to show on Friday.
'
$ws.Cells.Item(57, 3).Value = 'MSG: None

MSG: The rights to both movies have been acquired successfully.
'
$ws.Cells.Item(57, 4).Value = 'both_movies, '
$ws.Cells.Item(58, 3).Value = 'MSG: None

MSG: I have recorded the decision that no movie was selected for Friday.
'
$ws.Cells.Item(58, 4).Value = 'no_decision, '
